$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename destination-well-id values (shared strings 20/21): "_D" suffix -> "_E" suffix
$ws.Range("E2").Value = "ssdest000000141jul17_E"
$ws.Range("E3").Value = "ssdest000000141jul17_E"
$ws.Range("E4").Value = "ssdest000000141jul17_E"
$ws.Range("E5").Value = "ssdest000000141jul17_E"
$ws.Range("E6").Value = "ssdest000000141jul17_384_E"
$ws.Range("E7").Value = "ssdest000000141jul17_384_E"

# Add new "Dest Well Count" column (H) - clone the header formatting from G1
$ws.Range("H1").Value = "Dest Well Count"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H2").Value = 96
$ws.Range("H3").Value = 96
$ws.Range("H4").Value = 96
$ws.Range("H5").Value = 96
$ws.Range("H6").Value = 384
$ws.Range("H7").Value = 384

$ws.Columns.Item(8).ColumnWidth = 25.5

# Update selection to the new last cell
$null = $ws.Range("H7").Select()
